$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.403.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.504.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.31%  "
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0816"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.502.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.345.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.40%  "
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.00%  "
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "249.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.23%  "
$ws.Range("E30").Value = "  +6.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.139"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +1.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0799"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.35%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +6.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.008.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.60%  "
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  +5.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.57%  "
